# Read channel names dynamically from Summary sheet
#
# The fixture "Bad Report.xlsx" hardcoded a channel tab name ("Retail").
# Renaming it to "Unexpected" simulates a report with a channel name that
# isn't one of the hardcoded ones, so the parser must read channel names
# dynamically from the "Payments by Channel" section on the Summary sheet
# instead of assuming Retail/Subscription/Pool.
#
# This also reflects the author having the renamed sheet active/selected
# (and the previously active Summary sheet/cell no longer being the
# front-most selection) while preparing the fixture.

$wb = $excel.ActiveWorkbook

# Move the Summary sheet's selection off of its old spot before we hand
# focus to the renamed sheet.
$summary = $wb.Worksheets.Item("Summary")
[void]$summary.Range("A12").Select()

# Rename the "Retail" channel tab to "Unexpected" (sheetId / position stay
# put - only the visible name changes).
$channel = $wb.Worksheets.Item("Retail")
$channel.Name = "Unexpected"

# Make the renamed sheet the active tab, with its frozen-pane selection
# moved to N12.
$channel.Activate()
[void]$channel.Range("N12").Select()
